$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.325.40"
$ws.Range("E2").Value = "  -3.63%  "

$ws.Range("D3").Value = "3.156.67"
$ws.Range("E3").Value = "  -3.03%  "

$ws.Range("E4").Value = "  +0.11%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "606.98"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.16%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "147.50"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -6.78%  "

$ws.Range("E7").Value = "  +0.08%  "

$ws.Range("D8").Value = "3.154.28"
$ws.Range("E8").Value = "  -3.16%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.526"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -4.03%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.151"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -6.40%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.51"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -6.64%  "

$ws.Range("E12").Value = "  -5.96%  "

$ws.Range("E13").Value = "  -7.90%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "35.63"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -9.32%  "

$ws.Range("D15").Value = "3.679.04"
$ws.Range("E15").Value = "  -2.87%  "

$ws.Range("D16").Value = "64.319.41"
$ws.Range("E16").Value = "  -3.64%  "

$ws.Range("E17").Value = "  +0.30%  "

$ws.Range("D18").Value = "3.158.00"
$ws.Range("E18").Value = "  -2.78%  "

$ws.Range("E19").Value = "  -6.67%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "481.58"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -5.32%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.73"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -4.74%  "

$ws.Range("E22").Value = "  -5.65%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.78"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -4.33%  "

$ws.Range("E24").Value = "  -7.95%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "83.80"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -3.32%  "

$ws.Range("E26").Value = "  -0.03%  "

$ws.Range("E27").Value = "  -5.39%  "

$ws.Range("E28").Value = "  -7.51%  "

$ws.Range("E29").Value = "  -8.97%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.79"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.97%  "

$ws.Range("E31").Value = "  -20.51%  "

$ws.Range("E32").Value = "  -5.89%  "

$ws.Range("E33").Value = "  +0.07%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "26.28"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -6.52%  "

$ws.Range("E35").Value = "  -4.19%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "54.61"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.29%  "

$ws.Range("E37").Value = "  -7.44%  "

$ws.Range("D38").Value = "0.0₃0731"
$ws.Range("E38").Value = "  -8.64%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "452.12"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -8.96%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.94"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -12.94%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0396"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -7.48%  "

$ws.Range("E42").Value = "  -4.68%  "

$ws.Range("E43").Value = "  -8.56%  "

$ws.Range("D44").Value = "2.855.41"
$ws.Range("E44").Value = "  -3.92%  "

$ws.Range("E45").Value = "  -9.63%  "

$ws.Range("E46").Value = "  -9.60%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "26.58"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -8.10%  "

$ws.Range("E48").Value = "  +0.00%  "

$ws.Range("E49").Value = "  -7.30%  "

$ws.Range("E50").Value = "  -4.91%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "119.74"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.59%  "
